$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for all data rows (2-103)
# from serial date 45186 (2023-09-17) to 45188 (2023-09-19).
for ($row = 2; $row -le 103; $row++) {
    $ws.Cells.Item($row, 3).Value = 45188
}
